$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 19.56870467783557
$ws.Range("C2").Value = 0.02894430293274442
$ws.Range("B3").Value = 19.27198683555866
$ws.Range("C3").Value = 0.01149361143113127
